$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in column H, copying the header formatting (bold,
# border, centered) from the neighboring "sum" header cell (G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Add the new Save values for the two data rows (both 0, matching the diff)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
